$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$wsWeights = $wb.Worksheets.Item("network_weights")

# Row 1 had duplicate "value" header cells in C1:F1 left over from a copy/fill;
# trim the row back down to just A1:B1.
$ws.Range("C1:F1").ClearContents()

# The "Model" parameter row is renamed to "production_function".
$ws.Cells.Item(8, 1).Value = "production_function"

# Insert a new optimization parameter row right after it: "L_curve" = 1,
# matching the B8 "1/0 flag" number formatting used by its neighbours.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "L_curve"
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 2).NumberFormat = "0.00E+00"

# The obsolete "Deletion" row (now shifted down to row 17 by the insert above)
# is removed entirely.
$ws.Rows.Item(17).Delete()

# Reflect the saved selection/view state: row 17 (now simulation_timepoints)
# is selected on this sheet, and this sheet becomes the active tab instead of
# network_weights.
$ws.Range("A17:XFD17").Select()
$ws.Activate()
